$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New / updated data for rows 2-7 (columns A-T), following the revised NATMI LR-pair calc
$data = @(
  @("ECs", "Bmp6", "Bmpr1b", "FAPs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 1.57938, 4.73814, 0.7235635290775982, 0.7235635290775982, 50.26455976938, 452.38103792442, 0.438740623692882, 0.438740623692882),
  @("ECs", "Bmp6", "Bmpr1b", "sCs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 0.6034, 1.8102, 0.2764364709224018, 0.2764364709224018, 19.2035073034, 172.8315657306, 0.1676202638606827, 0.1676202638606827),
  @("FAPs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 1.57938, 4.73814, 0.7235635290775982, 0.7235635290775982, 7.654118759319998, 68.88706883388, 0.06680995225445434, 0.06680995225445434),
  @("FAPs", "Bmp6", "Bmpr1b", "sCs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 0.6034, 1.8102, 0.2764364709224018, 0.2764364709224018, 2.924245754266666, 26.3182117884, 0.02552465219917801, 0.02552465219917801),
  @("sCs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 1.57938, 4.73814, 0.7235635290775982, 0.7235635290775982, 24.9767733402, 224.7909600618, 0.2180129531302619, 0.218012953130262),
  @("sCs", "Bmp6", "Bmpr1b", "sCs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 0.6034, 1.8102, 0.2764364709224018, 0.2764364709224018, 9.542342586, 85.881083274, 0.08329155486254106, 0.08329155486254107)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}